$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking price values so Excel does not
# auto-convert them from text to numbers (the source data keeps them as text).
$textCells = @('D5','D6','D7','D8','D9','D10','D11','D12','D13','D14','D16','D17','D18','D19','D20','D21','D22','D23','D24','D26','D27','D28','D30','D31','D32','D33','D34','D35','D36','D37','D38','D39','D40','D41','D42','D43','D44','D45','D46','D47','D48','D49','D50','D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '20.543.84'
$ws.Range('E2').Value = '  +2.51%  '
$ws.Range('D3').Value = '1.468.08'
$ws.Range('E3').Value = '  +3.02%  '
$ws.Range('E4').Value = '  +0.79%  '
$ws.Range('D5').Value = '0.9515'
$ws.Range('E5').Value = '  -5.12%  '
$ws.Range('D6').Value = '281.52'
$ws.Range('E6').Value = '  +2.98%  '
$ws.Range('D7').Value = '0.3703'
$ws.Range('E7').Value = '  -1.19%  '
$ws.Range('D8').Value = '0.3187'
$ws.Range('E8').Value = '  +3.01%  '
$ws.Range('D9').Value = '41.86'
$ws.Range('E9').Value = '  +4.32%  '
$ws.Range('D10').Value = '1.059'
$ws.Range('D11').Value = '0.06681'
$ws.Range('E11').Value = '  +1.46%  '
$ws.Range('D12').Value = '1.005'
$ws.Range('E12').Value = '  +0.28%  '
$ws.Range('D13').Value = '5.629'
$ws.Range('E13').Value = '  +4.68%  '
$ws.Range('D14').Value = '18.21'
$ws.Range('E14').Value = '  +6.39%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '1.476.55'
$ws.Range('E15').Value = '  +3.22%  '
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').Value = '6.266'
$ws.Range('E16').Value = '  +1.64%  '
$ws.Range('D17').Value = '0.00001034'
$ws.Range('E17').Value = '  +2.16%  '
$ws.Range('D18').Value = '0.05670'
$ws.Range('E18').Value = '  -2.61%  '
$ws.Range('D19').Value = '72.23'
$ws.Range('E19').Value = '  -4.13%  '
$ws.Range('D20').Value = '0.9496'
$ws.Range('E20').Value = '  -5.26%  '
$ws.Range('D21').Value = '5.690'
$ws.Range('E21').Value = '  +0.36%  '
$ws.Range('D22').Value = '14.69'
$ws.Range('E22').Value = '  +1.63%  '
$ws.Range('D23').Value = '11.21'
$ws.Range('E23').Value = '  +1.62%  '
$ws.Range('D24').Value = '2.280'
$ws.Range('E24').Value = '  -2.23%  '
$ws.Range('D25').Value = '20.666.80'
$ws.Range('E25').Value = '  +3.07%  '
$ws.Range('D26').Value = '2.301'
$ws.Range('E26').Value = '  +0.72%  '
$ws.Range('D27').Value = '137.60'
$ws.Range('E27').Value = '  -0.87%  '
$ws.Range('D28').Value = '17.57'
$ws.Range('E28').Value = '  +4.28%  '
$ws.Range('D29').Value = '1.638.97'
$ws.Range('E29').Value = '  +3.11%  '
$ws.Range('D30').Value = '113.76'
$ws.Range('E30').Value = '  +4.49%  '
$ws.Range('D31').Value = '3.956'
$ws.Range('E31').Value = '  +2.68%  '
$ws.Range('D32').Value = '5.320'
$ws.Range('E32').Value = '  -1.54%  '
$ws.Range('D33').Value = '0.8395'
$ws.Range('E33').Value = '  -6.09%  '
$ws.Range('D34').Value = '1.635'
$ws.Range('E34').Value = '  +27.68%  '
$ws.Range('D35').Value = '0.07856'
$ws.Range('E35').Value = '  +1.07%  '
$ws.Range('D36').Value = '0.06065'
$ws.Range('E36').Value = '  +6.37%  '
$ws.Range('D37').Value = '4.921'
$ws.Range('E37').Value = '  +3.86%  '
$ws.Range('D38').Value = '10.63'
$ws.Range('E38').Value = '  -6.24%  '
$ws.Range('D39').Value = '0.02068'
$ws.Range('E39').Value = '  +2.42%  '
$ws.Range('D40').Value = '1.121'
$ws.Range('E40').Value = '  +1.06%  '
$ws.Range('D41').Value = '0.9625'
$ws.Range('E41').Value = '  -3.96%  '
$ws.Range('D42').Value = '0.1893'
$ws.Range('E42').Value = '  -1.13%  '
$ws.Range('D43').Value = '7.418'
$ws.Range('E43').Value = '  -11.74%  '
$ws.Range('D44').Value = '0.5424'
$ws.Range('E44').Value = '  +1.89%  '
$ws.Range('D45').Value = '12.59'
$ws.Range('E45').Value = '  +3.08%  '
$ws.Range('D46').Value = '3.589'
$ws.Range('E46').Value = '  +1.39%  '
$ws.Range('D47').Value = '122.39'
$ws.Range('E47').Value = '  +11.59%  '
$ws.Range('D48').Value = '0.5337'
$ws.Range('E48').Value = '  +4.07%  '
$ws.Range('D49').Value = '1.835'
$ws.Range('E49').Value = '  +3.09%  '
$ws.Range('D50').Value = '0.06439'
$ws.Range('E50').Value = '  +4.23%  '
$ws.Range('D51').Value = '1.053'
$ws.Range('E51').Value = '  +0.31%  '

# Restore default (General) formatting/style so the cells keep their
# original style (no explicit style index), matching the source formatting.
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "General"
    $ws.Range($addr).Style = "Normal"
}
